# Update countries & provincias Spain
# Refresh of the COVID-19 "Pais" stats sheet:
#  - the "last updated" timestamp moves from 04:26 to 05:43
#  - Birmania overtakes Eslovaquia in the case ranking (rows 111/112 swap identity)
#  - Islas Malvinas overtakes Montserrat in the case ranking (rows 214/215 swap identity)
#  - several countries get refreshed case/recovered/death counts

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text / shared-string updates ---
$ws.Range("A1").Value = "Datos actualizados a 23 de Septiembre de 2020 a las 05:43"

# Row 111/112: Eslovaquia / Birmania swap places (ranking overtake)
$ws.Range("A111").Value = "Birmania"
$ws.Range("A112").Value = "Eslovaquia"

# Row 214/215: Montserrat / Islas Malvinas swap places
$ws.Range("A214").Value = "Islas Malvinas"
$ws.Range("A215").Value = "Montserrat"

# --- Numeric updates ---
# Row 21
$ws.Range("B21").Value = 307418
$ws.Range("C21").Value = 532
$ws.Range("D21").Value = 293916
$ws.Range("E21").Value = 7070
$ws.Range("G21").Value = 8
$ws.Range("H21").Value = 6432

# Row 35
$ws.Range("B35").Value = 107450
$ws.Range("C35").Value = 76
$ws.Range("E35").Value = 3687
$ws.Range("H35").Value = 1699

# Row 38
$ws.Range("B38").Value = 105226
$ws.Range("C38").Value = 1834
$ws.Range("D38").Value = 19039
$ws.Range("E38").Value = 76232
$ws.Range("G38").Value = 5
$ws.Range("H38").Value = 9955

# Row 50
$ws.Range("B50").Value = 72306
$ws.Range("C50").Value = 231
$ws.Range("D50").Value = 23230
$ws.Range("E50").Value = 46870
$ws.Range("G50").Value = 2
$ws.Range("H50").Value = 2206

# Row 111
$ws.Range("B111").Value = 6959
$ws.Range("C111").Value = 216
$ws.Range("D111").Value = 1951
$ws.Range("E111").Value = 4892
$ws.Range("G111").Value = 1
$ws.Range("H111").Value = 116

# Row 112
$ws.Range("B112").Value = 6931
$ws.Range("D112").Value = 3668
$ws.Range("E112").Value = 3223
$ws.Range("H112").Value = 40

# Row 172
$ws.Range("B172").Value = 672
$ws.Range("C172").Value = 4
$ws.Range("D172").Value = 581
$ws.Range("E172").Value = 86

# Row 173
$ws.Range("B173").Value = 594
$ws.Range("C173").Value = 3
$ws.Range("D173").Value = 512
$ws.Range("E173").Value = 62

# Row 188
$ws.Range("D188").Value = 195
$ws.Range("E188").Value = 66

# Row 214
$ws.Range("D214").Value = 13
$ws.Range("H214").Value = 0

# Row 215
$ws.Range("D215").Value = 12
$ws.Range("H215").Value = 1

Write-Host "Update complete."
Write-Host ("A1: " + $ws.Range("A1").Value())
Write-Host ("A111 (" + $ws.Range("A111").Value() + "): " + $ws.Range("B111").Value() + ", " + $ws.Range("C111").Value() + ", " + $ws.Range("D111").Value() + ", " + $ws.Range("E111").Value() + ", " + $ws.Range("G111").Value() + ", " + $ws.Range("H111").Value())
Write-Host ("A112 (" + $ws.Range("A112").Value() + "): " + $ws.Range("B112").Value() + ", " + $ws.Range("C112").Value() + ", " + $ws.Range("D112").Value() + ", " + $ws.Range("E112").Value() + ", " + $ws.Range("G112").Value() + ", " + $ws.Range("H112").Value())
Write-Host ("A214 (" + $ws.Range("A214").Value() + "): " + $ws.Range("B214").Value() + ", " + $ws.Range("C214").Value() + ", " + $ws.Range("D214").Value() + ", " + $ws.Range("E214").Value() + ", " + $ws.Range("G214").Value() + ", " + $ws.Range("H214").Value())
Write-Host ("A215 (" + $ws.Range("A215").Value() + "): " + $ws.Range("B215").Value() + ", " + $ws.Range("C215").Value() + ", " + $ws.Range("D215").Value() + ", " + $ws.Range("E215").Value() + ", " + $ws.Range("G215").Value() + ", " + $ws.Range("H215").Value())

